$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column E (Payment type) values: row 3 -> "Remaining", all other data rows -> "Done"
for ($r = 2; $r -le 16; $r++) {
    if ($r -eq 3) {
        $ws.Cells.Item($r, 5).Value = "Remaining"
    } else {
        $ws.Cells.Item($r, 5).Value = "Done"
    }
}

# Make the formatting of E10:E16 match E2:E9 (copy formats down)
$ws.Range("E2").Copy()
$ws.Range("E10:E16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selection to E4:E16 with active cell E4
$ws.Range("E4:E16").Select()
